# Append a new record (row 84) to the "Optical_Power" sheet, mirroring the
# automated export/update process described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowIndex = 84

# Force the new row to "Text" format first so values that look like numbers
# or dates (e.g. "-506", "7/11/2025", "6", "808150511", "1") are stored as
# literal text, matching the rest of the sheet's data.
$rng = $ws.Range("A" + $rowIndex + ":P" + $rowIndex)
$rng.NumberFormat = "@"

$ws.Cells.Item($rowIndex, 1).Value  = "-506"
$ws.Cells.Item($rowIndex, 2).Value  = "7/11/2025"
$ws.Cells.Item($rowIndex, 3).Value  = "Gervasio Espinosa 591"
$ws.Cells.Item($rowIndex, 4).Value  = "6"
$ws.Cells.Item($rowIndex, 5).Value  = "808150511"
$ws.Cells.Item($rowIndex, 6).Value  = "Optical Power"
$ws.Cells.Item($rowIndex, 7).Value  = "Pendiente"
$ws.Cells.Item($rowIndex, 8).Value  = "Picada"
$ws.Cells.Item($rowIndex, 9).Value  = "1"
$ws.Cells.Item($rowIndex, 10).Value = "Cambio"
$ws.Cells.Item($rowIndex, 11).Value = "Nodo Teco"
$ws.Cells.Item($rowIndex, 12).Value = "Pasante"
$ws.Cells.Item($rowIndex, 13).Value = ""
$ws.Cells.Item($rowIndex, 14).Value = ""
$ws.Cells.Item($rowIndex, 15).Value = "No ubicado"
$ws.Cells.Item($rowIndex, 16).Value = "No clasificado, consultar con mantenimiento"

# Restore the default (unstyled) cell style on the new row so it matches the
# plain formatting used by the rest of the data rows.
$rng.Style = "Normal"
